$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bug report entry: "05.stats.tabs"
$ws.Range("A15").Value = "05.stats.tabs"
$ws.Range("D15").Value = "Tabs are not switching properly on reload and on clicking"
$ws.Range("D16").Value = "the stats navbar link"
$ws.Range("J15").Value = "home->statspage->tabs"
$ws.Range("O15").Value = "not required"
$ws.Range("T15").Value = "Active"

# Copy the "Active" status formatting (red fill / bold white font) from the
# existing row (T12) onto the new row's status cell (T15)
$ws.Range("T12").Copy()
$ws.Range("T15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the selection to where the user last left the cursor
$null = $ws.Range("Q20").Select()
